$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, shifting existing rows 23..147 down to 24..148.
$ws.Rows("23:23").Insert()

# Populate the newly inserted row 23 with a new weekly observation.
# Columns that stay identical to the row that used to occupy this slot
# (now shifted to row 24): Mercado ID, Mercado, Region, Codreg, Categoria ID,
# Categoria, Variedad, Calidad, Unidad de comercializacion, Origen,
# Kg o Unidades, Clasificacion.
$ws.Cells.Item(23, 1).Value = 10
$ws.Cells.Item(23, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(23, 3).Value = "La Araucanía"
$ws.Cells.Item(23, 4).Value = 44462
$ws.Cells.Item(23, 5).Value = 9
$ws.Cells.Item(23, 6).Value = 100112039
$ws.Cells.Item(23, 7).Value = "Ciboulette"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 40
$ws.Cells.Item(23, 11).Value = 6000
$ws.Cells.Item(23, 12).Value = 7000
$ws.Cells.Item(23, 13).Value = 6500
$ws.Cells.Item(23, 14).Value = "$/docena de atados"
$ws.Cells.Item(23, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(23, 16).Value = 2167
$ws.Cells.Item(23, 17).Value = 3
$ws.Cells.Item(23, 18).Value = "Hortaliza"
